$wb = $excel.ActiveWorkbook

$wsReroute = $wb.Worksheets.Item("Reroute Request")
$wsBol = $wb.Worksheets.Item("BOL")

# Update BOL order/tracking numbers and carrier code (manage order BOL fixes)
$wsBol.Range("A2").NumberFormat = "@"
$wsBol.Range("A2").Value = "51541518"

$wsBol.Range("A3").NumberFormat = "@"
$wsBol.Range("A3").Value = "51541536"

$wsBol.Range("A4").NumberFormat = "@"
$wsBol.Range("A4").Value = "51541536"

$wsBol.Range("C3").Value = "cev"
$wsBol.Range("C4").Value = "cev"

# Update sheet selections
$wsReroute.Range("A4").Select() | Out-Null
$wsBol.Range("B4").Select() | Out-Null

# BOL sheet becomes the active/selected tab
$wsBol.Activate() | Out-Null
